$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$title = $s.Shapes.Item(1)
# TextRange.Text already reads back as the concatenation "Two-Column Layout"
# (it spans 3 runs: "Two-Column", " ", "Layout"), so assigning that exact
# same string directly is a no-op and the run split survives. Stage a
# throwaway value first (distinct, non-empty, so the paragraph keeps a run
# to inherit default formatting from instead of stamping a fresh lang
# attribute) so the following assignment performs a real text replace and
# collapses the paragraph down to a single run.
$title.TextFrame.TextRange.Text = "placeholder"
$title.TextFrame.TextRange.Text = "Two-Column Layout"
